$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.433768000000001
$ws.Range("H2").Value = 16.301304
$ws.Range("I2").Value = 0.1262505823713576
$ws.Range("J2").Value = 0.1262505823713576
$ws.Range("M2").Value = 7.655977
$ws.Range("N2").Value = 22.967931
$ws.Range("O2").Value = 0.2994795900616967
$ws.Range("P2").Value = 0.2994795900616967
$ws.Range("Q2").Value = 41.600802831336
$ws.Range("R2").Value = 374.407225482024
$ws.Range("S2").Value = 0.03780947265362465
$ws.Range("T2").Value = 0.03780947265362464

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.433768000000001
$ws.Range("H3").Value = 16.301304
$ws.Range("I3").Value = 0.1262505823713576
$ws.Range("J3").Value = 0.1262505823713576
$ws.Range("O3").Value = 0.3140620915319453
$ws.Range("P3").Value = 0.3140620915319453
$ws.Range("Q3").Value = 43.62646263782401
$ws.Range("R3").Value = 392.638163740416
$ws.Range("S3").Value = 0.03965052195667473
$ws.Range("T3").Value = 0.03965052195667471

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.433768000000001
$ws.Range("H4").Value = 16.301304
$ws.Range("I4").Value = 0.1262505823713576
$ws.Range("J4").Value = 0.1262505823713576
$ws.Range("M4").Value = 9.879524666666667
$ws.Range("N4").Value = 29.638574
$ws.Range("O4").Value = 0.386458318406358
$ws.Range("P4").Value = 0.386458318406358
$ws.Range("Q4").Value = 53.683044988944
$ws.Range("R4").Value = 483.147404900496
$ws.Range("S4").Value = 0.04879058776105826
$ws.Range("T4").Value = 0.04879058776105825

$ws.Range("G5").Value = 5.270503666666666
$ws.Range("I5").Value = 0.1224572262391479
$ws.Range("J5").Value = 0.1224572262391479
$ws.Range("M5").Value = 7.655977
$ws.Range("N5").Value = 22.967931
$ws.Range("O5").Value = 0.2994795900616967
$ws.Range("P5").Value = 0.2994795900616967
$ws.Range("Q5").Value = 40.35085485041566
$ws.Range("R5").Value = 363.157693653741
$ws.Range("S5").Value = 0.03667343991419247
$ws.Range("T5").Value = 0.03667343991419246

$ws.Range("G6").Value = 5.270503666666666
$ws.Range("I6").Value = 0.1224572262391479
$ws.Range("J6").Value = 0.1224572262391479
$ws.Range("O6").Value = 0.3140620915319453
$ws.Range("P6").Value = 0.3140620915319453
$ws.Range("S6").Value = 0.03845917259586741
$ws.Range("T6").Value = 0.0384591725958674

$ws.Range("G7").Value = 5.270503666666666
$ws.Range("I7").Value = 0.1224572262391479
$ws.Range("J7").Value = 0.1224572262391479
$ws.Range("M7").Value = 9.879524666666667
$ws.Range("N7").Value = 29.638574
$ws.Range("O7").Value = 0.386458318406358
$ws.Range("P7").Value = 0.386458318406358
$ws.Range("Q7").Value = 52.07007098059044
$ws.Range("R7").Value = 468.630638825314
$ws.Range("S7").Value = 0.04732461372908805
$ws.Range("T7").Value = 0.04732461372908804

$ws.Range("G8").Value = 32.33527633333333
$ws.Range("H8").Value = 97.00582900000001
$ws.Range("I8").Value = 0.7512921913894945
$ws.Range("J8").Value = 0.7512921913894944
$ws.Range("M8").Value = 7.655977
$ws.Range("N8").Value = 22.967931
$ws.Range("O8").Value = 0.2994795900616967
$ws.Range("P8").Value = 0.2994795900616967
$ws.Range("Q8").Value = 247.5581318966443
$ws.Range("R8").Value = 2228.023187069799
$ws.Range("S8").Value = 0.2249966774938796
$ws.Range("T8").Value = 0.2249966774938796

$ws.Range("G9").Value = 32.33527633333333
$ws.Range("H9").Value = 97.00582900000001
$ws.Range("I9").Value = 0.7512921913894945
$ws.Range("J9").Value = 0.7512921913894944
$ws.Range("O9").Value = 0.3140620915319453
$ws.Range("P9").Value = 0.3140620915319453
$ws.Range("Q9").Value = 259.612431896224
$ws.Range("R9").Value = 2336.511887066016
$ws.Range("S9").Value = 0.2359523969794032
$ws.Range("T9").Value = 0.2359523969794032

$ws.Range("G10").Value = 32.33527633333333
$ws.Range("H10").Value = 97.00582900000001
$ws.Range("I10").Value = 0.7512921913894945
$ws.Range("J10").Value = 0.7512921913894944
$ws.Range("M10").Value = 9.879524666666667
$ws.Range("N10").Value = 29.638574
$ws.Range("O10").Value = 0.386458318406358
$ws.Range("P10").Value = 0.386458318406358
$ws.Range("Q10").Value = 319.4571601386496
$ws.Range("R10").Value = 2875.114441247846
$ws.Range("S10").Value = 0.2903431169162117
$ws.Range("T10").Value = 0.2903431169162117
